$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Find the last used row so we know how far to search.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# Collect the row numbers of every "Scenario type" setting (one per scenario
# block) by scanning column B. We search from the bottom up so that
# inserting rows for a later scenario does not invalidate the row numbers
# we already found for earlier scenarios.
$scenarioTypeRows = New-Object System.Collections.ArrayList
for ($r = 1; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Value2
    if ($label -eq "Scenario type") {
        [void]$scenarioTypeRows.Add($r)
    }
}

for ($i = $scenarioTypeRows.Count - 1; $i -ge 0; $i--) {
    $r = $scenarioTypeRows[$i]
    $insertAt = $r + 1

    # Insert two new blank rows right after the "Scenario type" row, and
    # before the following "Substrate" row.
    $ws.Range("A" + $insertAt + ":A" + ($insertAt + 1)).EntireRow.Insert()

    $ws.Cells.Item($insertAt, 1).Value = 1
    $ws.Cells.Item($insertAt, 2).Value = "Slice thickness"
    $ws.Cells.Item($insertAt, 3).Value = 0

    $ws.Cells.Item($insertAt + 1, 1).Value = 1
    $ws.Cells.Item($insertAt + 1, 2).Value = "Slice thickness unit"
    $ws.Cells.Item($insertAt + 1, 3).Value = "um"
}
